$wb = $excel.ActiveWorkbook

# --- Update "Status" text from "Handed back: in sync with en-US" to "Ready for handoff" ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Update Latest Handoff / Xliff Generate datetimes ---
# Overview: Latest HO Xliff Generate Date (G2) 2017-02-21 05:05:06 -> 2017-02-21 05:07:18
$wsOverview.Range("G2").Value = "2017-02-21 05:07:18"

# de-de: Latest Handoff Datetime (H2) 2017-02-21 05:05:06 -> 2017-02-21 05:07:18
$wsDeDe.Range("H2").Value = "2017-02-21 05:07:18"

# zh-cn: Latest Handoff Datetime (H2) 2017-02-21 05:04:51 -> 2017-02-21 05:07:02
$wsZhCn.Range("H2").Value = "2017-02-21 05:07:02"

# --- Column widths shrink to fit the shorter status text ---
# (target stored width is 17.2159881591797; the host snaps ColumnWidth to a
# 1/6-character grid, so 16.33 is the input that lands closest to it)
$wsOverview.Columns.Item(5).ColumnWidth = 16.33
$wsOverview.Columns.Item(6).ColumnWidth = 16.33

$wsZhCn.Columns.Item(3).ColumnWidth = 16.33

$wsDeDe.Columns.Item(3).ColumnWidth = 16.33

$wb.Save()
